$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string / text edits -------------------------------------------------
# C5 first (brand-new cell) gets the "Remove previous marker" text -> appended as
# a new shared string; then C3 (existing cell, sole referrer of the old string)
# is retargeted to the typo-fixed text -> frees/repacks the old slot and appends
# the typo-fixed text last, landing the two strings in the exact order the
# original author's workbook has them in.
$ws.Range("C5").Value2 = "Fix bug: Remove previous marker when another accident spot is clicked."
$ws.Range("C3").Value2 = "Fix bug: The info window of previous marker doesn't close automatically when another marker is clicked."

# --- New row 3 date (A3) ----------------------------------------------------
$ws.Range("A3").Value2 = 43895

# --- Row 2 extra column (HOURS) ---------------------------------------------
$ws.Range("F2").Value2 = 0.5

# --- Row 3 remaining values ---------------------------------------------------
$ws.Range("F3").Value2 = 0.25

# --- Row 5 (new) --------------------------------------------------------------
$ws.Range("A5").Value2 = 43899
$ws.Range("B5").Value2 = "1"
$ws.Range("D5").Value2 = 0.52083333333333337
$ws.Range("E5").Value2 = 0.53125
$ws.Range("F5").Value2 = 0.25

# --- Row 6 (new) --------------------------------------------------------------
$ws.Range("B6").Value2 = "2"
$ws.Range("C6").Value2 = "New feature: Show nearest k fire stations"
$ws.Range("D6").Value2 = 0.53125
$ws.Range("E6").Value2 = 0.57291666666666663
$ws.Range("F6").Value2 = 1

# --- Row 7 (new) --------------------------------------------------------------
$ws.Range("B7").Value2 = "3"
$ws.Range("C7").Value2 = "New feature: Add custom control"
$ws.Range("D7").Value2 = 0.57291666666666663

# --- Row 8 (new) --------------------------------------------------------------
$ws.Range("B8").Value2 = "4"
$ws.Range("C8").Value2 = "Fix bug: Restrict the map bounders"
$ws.Range("D8").Value2 = 0.58333333333333337
$ws.Range("E8").Value2 = 0.59375
$ws.Range("F8").Value2 = 0.25
